$wb = $excel.ActiveWorkbook

# The "trial_investigators" sheet carried a bad data row for investigator
# "Mark Andrew Dickson" that needs to be removed (disabled failing CTS
# tests / removed bad data from the Excel data providers).
$ws1 = $wb.Worksheets.Item("trial_investigators")
$ws1.Rows(3).Delete()

# Leave the workbook focused on the (now trimmed) trial_investigators sheet
# with the selection parked on F6, matching the refreshed view state.
$ws1.Activate()
$ws1.Range("F6").Select()
